# Scheduled-runner market data refresh: update current average price /
# leve price / leve profit figures per sheet ("Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 206.66667
$ws.Range("I41").Value = 158
$ws.Range("J41").Value = 450
$ws.Range("K41").Value = 158
$ws.Range("L41").Value = 450
$ws.Range("M41").Value = 282
$ws.Range("N41").Value = -1330

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2666.6667
$ws.Range("J43").Value = 3000
$ws.Range("L43").Value = 3000
$ws.Range("N43").Value = -3138

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3300

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3300

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7999.8335
$ws.Range("I76").Value = 7999
$ws.Range("J76").Value = 8000
$ws.Range("K76").Value = 7999
$ws.Range("L76").Value = 8000
$ws.Range("M76").Value = -7684
$ws.Range("N76").Value = -8630

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 7999.8335
$ws.Range("I79").Value = 7999
$ws.Range("J79").Value = 8000
$ws.Range("K79").Value = 7999
$ws.Range("L79").Value = 8000
$ws.Range("M79").Value = -6907
$ws.Range("N79").Value = -10184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 498.85715
$ws.Range("I107").Value = 518
$ws.Range("K107").Value = 518
$ws.Range("M107").Value = 1402

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2461
$ws.Range("I112").Value = 1707.5
$ws.Range("K112").Value = 5122.5
$ws.Range("M112").Value = -4014.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5092.7
$ws.Range("I116").Value = 3419.7144
$ws.Range("K116").Value = 3419.7144
$ws.Range("M116").Value = 22.28560000000016

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5345.092
$ws.Range("I138").Value = 7115.421
$ws.Range("J138").Value = 4850.4414
$ws.Range("K138").Value = 21346.263
$ws.Range("L138").Value = 14551.3242
$ws.Range("M138").Value = -16206.263
$ws.Range("N138").Value = -24831.3242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16219.355
$ws.Range("I32").Value = 6598.2
$ws.Range("K32").Value = 6598.2
$ws.Range("M32").Value = -6311.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2227.875
$ws.Range("J45").Value = 3000
$ws.Range("L45").Value = 3000
$ws.Range("N45").Value = -3754

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("N56").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3500
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 5772.8823
$ws.Range("I110").Value = 5364.6924
$ws.Range("K110").Value = 5364.6924
$ws.Range("M110").Value = -3319.6924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 503930.56
$ws.Range("I122").Value = 836801
$ws.Range("K122").Value = 2510403
$ws.Range("M122").Value = -2507953

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3500
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3002.5745
$ws.Range("I105").Value = 2366.3235
$ws.Range("J105").Value = 4666.615
$ws.Range("K105").Value = 2366.3235
$ws.Range("L105").Value = 4666.615
$ws.Range("M105").Value = -619.3235
$ws.Range("N105").Value = -8160.615

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("N63").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("N66").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1100.875
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1100.875
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = ""
$ws.Range("N94").Value = -2002.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4990.3335
$ws.Range("I105").Value = 3827.3333
$ws.Range("K105").Value = 3827.3333
$ws.Range("M105").Value = -2080.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 852.4286
$ws.Range("I107").Value = 396.8889
$ws.Range("K107").Value = 396.8889
$ws.Range("M107").Value = 1523.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1999.3939
$ws.Range("I132").Value = 1805.8387
$ws.Range("K132").Value = 5417.5161
$ws.Range("M132").Value = -2887.5161

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2608.48
$ws.Range("I134").Value = 2198.9443
$ws.Range("K134").Value = 6596.8329
$ws.Range("M134").Value = -4061.8329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 848.75
$ws.Range("J52").Value = 848.75
$ws.Range("L52").Value = 2546.25
$ws.Range("N52").Value = -3078.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1567.3684
$ws.Range("J131").Value = 1634.1765
$ws.Range("L131").Value = 4902.529500000001
$ws.Range("N131").Value = -14982.5295

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5355.25
$ws.Range("I80").Value = 3750
$ws.Range("K80").Value = 3750
$ws.Range("M80").Value = -2752

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 5355.25
$ws.Range("I83").Value = 3750
$ws.Range("K83").Value = 18750
$ws.Range("M83").Value = -13758

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1728.1538
$ws.Range("I97").Value = 1870.75
$ws.Range("K97").Value = 1870.75
$ws.Range("M97").Value = -1374.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4996.636
$ws.Range("I126").Value = 4990.75
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 14972.25
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -12502.25
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("N138").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1873.6154
$ws.Range("I16").Value = 1896.091
$ws.Range("K16").Value = 1896.091
$ws.Range("M16").Value = -1726.091

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4358.5454
$ws.Range("I61").Value = 4099.4736
$ws.Range("K61").Value = 4099.4736
$ws.Range("M61").Value = -3897.4736

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4358.5454
$ws.Range("I113").Value = 4099.4736
$ws.Range("K113").Value = 4099.4736
$ws.Range("M113").Value = -1929.4736

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8153
$ws.Range("I122").Value = 6538.6
$ws.Range("K122").Value = 19615.8
$ws.Range("M122").Value = -17165.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4455.6895
$ws.Range("I132").Value = 3357.25
$ws.Range("K132").Value = 10071.75
$ws.Range("M132").Value = -7541.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3943.2856
$ws.Range("I136").Value = 3943.2856
$ws.Range("K136").Value = 11829.8568
$ws.Range("M136").Value = -9279.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1999.6666
$ws.Range("I113").Value = 1766
$ws.Range("J113").Value = 2233.3333
$ws.Range("K113").Value = 5298
$ws.Range("L113").Value = 6699.999899999999
$ws.Range("M113").Value = -3128
$ws.Range("N113").Value = -11039.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 97097
$ws.Range("I126").Value = 203917.4
$ws.Range("J126").Value = 8080
$ws.Range("K126").Value = 611752.2
$ws.Range("L126").Value = 24240
$ws.Range("M126").Value = -609282.2
$ws.Range("N126").Value = -29180
